$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "conversion" worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "conversion"

# Header row
$ws2.Range("A1").Value = "material_IRI"
$ws2.Range("B1").Value = "unit"
$ws2.Range("C1").Value = "input"
$ws2.Range("D1").Value = "unit"
$ws2.Range("E1").Value = "amout_input"
$ws2.Range("F1").Value = "input_iri"

# Row 2 - Cement
$ws2.Range("A2").Value = "http://data.europa.eu/ehl/cpa21/235"
$ws2.Range("B2").Value = "m3"
$ws2.Range("C2").Value = "Cement"
$ws2.Range("D2").Value = "kg"
$ws2.Range("E2").Value = 312
$ws2.Range("F2").Value = "http://data.europa.eu/ehl/cpa21/235"

# Row 3 - Gravel
$ws2.Range("A3").Value = "http://data.europa.eu/ehl/cpa21/235"
$ws2.Range("B3").Value = "m3"
$ws2.Range("C3").Value = "Gravel"
$ws2.Range("D3").Value = "kg"
$ws2.Range("E3").Value = 950
$ws2.Range("F3").Value = "http://data.europa.eu/ehl/cpa21/081"

# Row 4 - Sand
$ws2.Range("A4").Value = "http://data.europa.eu/ehl/cpa21/235"
$ws2.Range("B4").Value = "m3"
$ws2.Range("C4").Value = "Sand"
$ws2.Range("D4").Value = "kg"
$ws2.Range("E4").Value = 815
$ws2.Range("F4").Value = "http://data.europa.eu/ehl/cpa21/081"

# Match the final selection recorded in the source file
$ws2.Range("B8").Select()
